# ============================================================================
# Modulo de edicao de dados adicionado a faculdade
# - Adds a row to "Contas de Casa"
# - Fixes (buggy, as per the real commit) a row in "Anotacao Contas":
#   salvarDadosEditados swapped Mes/Ano into B/C as text and saved Valor as text
# - Removes the last (testing) row from "Anotacao Contas"
# - Adds a brand-new "Faculdade" worksheet with grade-tracking data
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet "Contas de Casa": append a new row (row 3)
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Contas de Casa")

$ws1.Cells.Item(3, 1).Value = "3e8f8227-8adf-4920-b55d-7a004e81ad98"
$ws1.Cells.Item(3, 2).Value = "teste"

$ws1.Cells.Item(3, 3).NumberFormat = "@"
$ws1.Cells.Item(3, 3).Value = "111"

$ws1.Cells.Item(3, 4).NumberFormat = "@"
$ws1.Cells.Item(3, 4).Value = "02/01/2024"

$ws1.Cells.Item(3, 5).NumberFormat = "@"
$ws1.Cells.Item(3, 5).Value = "01/01/2024"

$ws1.Cells.Item(3, 6).Value = "Sim"
$ws1.Cells.Item(3, 7).Value = "b"

# ----------------------------------------------------------------------------
# Sheet "Anotacao Contas": fix row 4 (Mes/Ano got swapped & stringified,
# Valor got stringified) and drop row 5 entirely
# ----------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Anotacao Contas")

$ws2.Cells.Item(4, 2).NumberFormat = "@"
$ws2.Cells.Item(4, 2).Value = "2024"

$ws2.Cells.Item(4, 3).Value = "Janeiro"

$ws2.Cells.Item(4, 6).NumberFormat = "@"
$ws2.Cells.Item(4, 6).Value = "521"

$ws2.Rows.Item(5).Delete()

# ----------------------------------------------------------------------------
# New sheet "Faculdade" (added after "Anotacao Contas", i.e. last tab)
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Faculdade"

$headers = @(
    "ID",
    "Nome da Materia",
    "Nota Atividade 1",
    "Nota Atividade 2",
    "Nota Atividade 3",
    "Nota Atividade 4",
    "Nota Mapa",
    "Nota SGC",
    "Valor Mensalidade",
    "Data Mensalidade",
    "Pago",
    "Unnamed: 11",
    "Unnamed: 12",
    "Unnamed: 13",
    "Unnamed: 14",
    "Nome Matéria",
    "Nota MAPA"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2
$ws3.Cells.Item(2, 1).Value = "7ad0e37b-302c-455f-aa93-eaf7d8228f33"
$ws3.Cells.Item(2, 2).Value = "Engenharia de Software LL"
$ws3.Cells.Item(2, 3).Value = 0.5
$ws3.Cells.Item(2, 4).Value = 0.5
$ws3.Cells.Item(2, 5).Value = 0.4
$ws3.Cells.Item(2, 6).Value = 0.3
$ws3.Cells.Item(2, 7).Value = 3.5
$ws3.Cells.Item(2, 8).Value = 0.1
$ws3.Cells.Item(2, 9).Value = 127

$ws3.Cells.Item(2, 10).NumberFormat = "@"
$ws3.Cells.Item(2, 10).Value = "01/01/2024"

$ws3.Cells.Item(2, 11).Value = "Sim"

# L2:O2 are blank placeholder numeric cells in the source data (NaN columns
# written out as 0 by the original export) - keep them as explicit zeros
$ws3.Cells.Item(2, 12).Value = 0
$ws3.Cells.Item(2, 13).Value = 0
$ws3.Cells.Item(2, 14).Value = 0
$ws3.Cells.Item(2, 15).Value = 0

$ws3.Cells.Item(2, 16).Value = "Engenharia de Software xx"
$ws3.Cells.Item(2, 17).Value = 3.5

# Row 3
$ws3.Cells.Item(3, 1).Value = "5667db64-9a3e-4d25-a4dd-df433cea8627"
$ws3.Cells.Item(3, 2).Value = "Engenharia de Software"
$ws3.Cells.Item(3, 3).Value = 0.5
$ws3.Cells.Item(3, 4).Value = 0.5
$ws3.Cells.Item(3, 5).Value = 0.4
$ws3.Cells.Item(3, 6).Value = 0.3
$ws3.Cells.Item(3, 7).Value = 3.5
$ws3.Cells.Item(3, 8).Value = 0.1
$ws3.Cells.Item(3, 9).Value = 127

$ws3.Cells.Item(3, 10).NumberFormat = "@"
$ws3.Cells.Item(3, 10).Value = "01/01/2024"

$ws3.Cells.Item(3, 11).Value = "Não"

$ws3.Cells.Item(3, 12).Value = 0
$ws3.Cells.Item(3, 13).Value = 0
$ws3.Cells.Item(3, 14).Value = 0
$ws3.Cells.Item(3, 15).Value = 0
$ws3.Cells.Item(3, 16).Value = 0
$ws3.Cells.Item(3, 17).Value = 0

# Header row formatting: bold, thin border all around, centered / top aligned
$headerRange = $ws3.Range("A1:Q1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
